{"js": "// Replace the 25 \"AxB=C\" answer cells in the multiplication-practice table\n// with their updated values. Each \"before\" string is unique in the\n// document, so a plain text search+replace (matchCase, no wildcards) is\n// unambiguous for every cell.\nconst replacements = [\n  [\"67\u00d753=3551\", \"44\u00d736=1584\"],\n  [\"40\u00d772=2880\", \"42\u00d772=3024\"],\n  [\"60\u00d753=3180\", \"23\u00d720=460\"],\n  [\"91\u00d753=4823\", \"28\u00d783=2324\"],\n  [\"95\u00d798=9310\", \"79\u00d768=5372\"],\n  [\"84\u00d736=3024\", \"88\u00d718=1584\"],\n  [\"65\u00d785=5525\", \"33\u00d788=2904\"],\n  [\"40\u00d755=2200\", \"54\u00d797=5238\"],\n  [\"75\u00d743=3225\", \"14\u00d721=294\"],\n  [\"79\u00d777=6083\", \"50\u00d791=4550\"],\n  [\"25\u00d785=2125\", \"97\u00d725=2425\"],\n  [\"89\u00d784=7476\", \"76\u00d798=7448\"],\n  [\"12\u00d785=1020\", \"65\u00d746=2990\"],\n  [\"18\u00d788=1584\", \"89\u00d726=2314\"],\n  [\"17\u00d791=1547\", \"74\u00d772=5328\"],\n  [\"24\u00d736=864\", \"65\u00d754=3510\"],\n  [\"94\u00d725=2350\", \"16\u00d781=1296\"],\n  [\"62\u00d771=4402\", \"98\u00d779=7742\"],\n  [\"39\u00d726=1014\", \"81\u00d767=5427\"],\n  [\"17\u00d741=697\", \"81\u00d794=7614\"],\n  [\"24\u00d755=1320\", \"11\u00d732=352\"],\n  [\"61\u00d711=671\", \"21\u00d726=546\"],\n  [\"34\u00d781=2754\", \"99\u00d768=6732\"],\n  [\"64\u00d769=4416\", \"85\u00d788=7480\"],\n  [\"78\u00d750=3900\", \"64\u00d785=5440\"],\n];\n\nconst body = context.document.body;\n\nfor (const [before, after] of replacements) {\n  const found = body.search(before, { matchCase: true, matchWholeWord: false });\n  found.load(\"items\");\n  await context.sync();\n\n  if (found.items.length === 0) {\n    throw new Error(`Text not found: ${before}`);\n  }\n\n  for (const range of found.items) {\n    range.insertText(after, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Replace the 25 \"AxB=C\" answer cells in the multiplication-practice table\n# with their updated values. Each \"before\" string is unique in the\n# document, so Find/Replace (MatchCase, whole document) is unambiguous for\n# every cell.\n$d = $word.ActiveDocument\n\n$pairs = @(\n    ,@(\"67\u00d753=3551\", \"44\u00d736=1584\")\n    ,@(\"40\u00d772=2880\", \"42\u00d772=3024\")\n    ,@(\"60\u00d753=3180\", \"23\u00d720=460\")\n    ,@(\"91\u00d753=4823\", \"28\u00d783=2324\")\n    ,@(\"95\u00d798=9310\", \"79\u00d768=5372\")\n    ,@(\"84\u00d736=3024\", \"88\u00d718=1584\")\n    ,@(\"65\u00d785=5525\", \"33\u00d788=2904\")\n    ,@(\"40\u00d755=2200\", \"54\u00d797=5238\")\n    ,@(\"75\u00d743=3225\", \"14\u00d721=294\")\n    ,@(\"79\u00d777=6083\", \"50\u00d791=4550\")\n    ,@(\"25\u00d785=2125\", \"97\u00d725=2425\")\n    ,@(\"89\u00d784=7476\", \"76\u00d798=7448\")\n    ,@(\"12\u00d785=1020\", \"65\u00d746=2990\")\n    ,@(\"18\u00d788=1584\", \"89\u00d726=2314\")\n    ,@(\"17\u00d791=1547\", \"74\u00d772=5328\")\n    ,@(\"24\u00d736=864\", \"65\u00d754=3510\")\n    ,@(\"94\u00d725=2350\", \"16\u00d781=1296\")\n    ,@(\"62\u00d771=4402\", \"98\u00d779=7742\")\n    ,@(\"39\u00d726=1014\", \"81\u00d767=5427\")\n    ,@(\"17\u00d741=697\", \"81\u00d794=7614\")\n    ,@(\"24\u00d755=1320\", \"11\u00d732=352\")\n    ,@(\"61\u00d711=671\", \"21\u00d726=546\")\n    ,@(\"34\u00d781=2754\", \"99\u00d768=6732\")\n    ,@(\"64\u00d769=4416\", \"85\u00d788=7480\")\n    ,@(\"78\u00d750=3900\", \"64\u00d785=5440\")\n)\n\nforeach ($pair in $pairs) {\n    $before = $pair[0]\n    $after = $pair[1]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $before\n    $find.Replacement.Text = $after\n    $find.MatchCase = $true\n    $find.MatchWholeWord = $false\n    $find.MatchWildcards = $false\n\n    $ok = $find.Execute($find.Text, $false, $false, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2)\n    if (-not $ok) {\n        throw \"Replace failed for: $before\"\n    }\n}\n"}
